$d = $word.ActiveDocument

# 1) Replace the bookmark-style placeholder text and drop the trailing
#    space run that followed it (the two runs share identical formatting,
#    so a single Find/Replace spanning both collapses them into one run).
$d.Content.Find.Execute("**ID__AFFARS_pgi_5343_topic_2__ID** ", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AF_PGI_5343__ID**", 2)

# 2) Add a paragraph border (top/left/bottom/right, 5-twip space, no
#    explicit line style) and tighten the left indent to 225 twips
#    (11.25 pt) on the first paragraph.
$p = $d.Paragraphs(1)
$p.Borders.DistanceFromTop = 5
$p.Borders.DistanceFromLeft = 5
$p.Borders.DistanceFromBottom = 5
$p.Borders.DistanceFromRight = 5
$p.LeftIndent = 11.25
